# defaultDialog.csv sheet: add food/sailor check dialog rows, and start
# wiring up the "hire sailors" feature.
#
# - widen column A (dialogId) so the longer new ids are readable
# - append 3 new dialog rows (no_enough_sailors, checkout_food,
#   checkout_food_not_enough) reusing the existing "name_dock_sailor"
#   dialogName and the npcType/npcParameter pairing already used by the
#   other dock_sailor rows (2, 7)
# - move the active selection to B7, where work continues next

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width 23.6640625 -> 31 (stored character-width units). The
# ColumnWidth COM property is expressed in points and rounds through
# Excel's MDW pixel formula, so 30.140625 is the value that lands on the
# exact stored width of 31.
$ws.Columns.Item(1).ColumnWidth = 30.140625

# Row 14: dialog_no_enough_sailors
$ws.Range("A14").Value = "dialog_no_enough_sailors"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "name_dock_sailor"
$ws.Range("F14").Value = 0

# Row 15: dialog_checkout_food
$ws.Range("A15").Value = "dialog_checkout_food"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = "name_dock_sailor"
$ws.Range("F15").Value = 0

# Row 16: dialog_checkout_food_not_enough
$ws.Range("A16").Value = "dialog_checkout_food_not_enough"
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = "name_dock_sailor"
$ws.Range("F16").Value = 0

# Leave the selection where the author left off editing.
$ws.Range("B7").Select()
